# ROBE-139 Enum type parsers
# Add a new "H" column to the sample sheet holding an enum-style sample
# value (SAMPLE1 / SAMPLE2, alternating) for each of the 5 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("SAMPLE1", "SAMPLE2", "SAMPLE1", "SAMPLE2", "SAMPLE1")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 8).Value = $values[$i]
}

# Leave the selection on the last-written cell, matching the recorded
# end-user state after entering the new column's data.
$ws.Range("H4").Select()
